# CIERRE 1 SEPT 23
# Advance the payroll receipt sheet from "SEMANA 34" (week of 21-27 Aug 2023)
# to "SEMANA 35" (week of 28 Aug - 03 Sep 2023): update the week-label cell,
# bump the first employee's EXTRAS amount, and zero out the third employee's
# EXTRAS amount. Formulas that reference these cells (sums, copies of the
# week label, and the TODAY()-based dates) recalculate on their own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("recibos")

# Week label, e.g. "SEMANA  34 ... AGOSTO 2023" -> "SEMANA  35 ... SEPTIEMBRE 2023"
$ws.Range("B9").Value = "SEMANA  35        DEL    28     Al   03  SEPTIEMBRE    2023"

# ALEJANDRA BAUTISTA SALAZAR's EXTRAS (K4): 280 -> 420
$ws.Range("K4").Value = 420

# TEODORA ARELLANO PEREZ's EXTRAS (E25): 933 -> 0
$ws.Range("E25").Value = 0

# Restore the on-screen scroll position to where row 37 is at the top.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 37
$ws.Range("E26").Select()
